# Update countries & provincias Spain
# - refresh the "last updated" timestamp
# - update COVID case counters for a number of countries
# - a few countries change rank (their row's displayed name changes
#   because the underlying data got re-sorted), so both the name (col A)
#   and the numbers (cols B-H) are rewritten for those rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 2 de Julio de 2020 a las 02:45"
$ws.Cells.Item(4, 2).Value = 2778500
$ws.Cells.Item(4, 3).Value = 50647
$ws.Cells.Item(4, 4).Value = 1159838
$ws.Cells.Item(4, 5).Value = 1487873
$ws.Cells.Item(4, 7).Value = 667
$ws.Cells.Item(4, 8).Value = 130789
$ws.Cells.Item(7, 2).Value = 605220
$ws.Cells.Item(7, 3).Value = 19428
$ws.Cells.Item(7, 4).Value = 359896
$ws.Cells.Item(7, 5).Value = 227476
$ws.Cells.Item(17, 2).Value = 196324
$ws.Cells.Item(17, 3).Value = 492
$ws.Cells.Item(17, 5).Value = 7463
$ws.Cells.Item(22, 4).Value = 67744
$ws.Cells.Item(22, 5).Value = 27912
$ws.Cells.Item(44, 2).Value = 34463
$ws.Cells.Item(44, 3).Value = 913
$ws.Cells.Item(44, 4).Value = 15945
$ws.Cells.Item(44, 5).Value = 17873
$ws.Cells.Item(44, 7).Value = 14
$ws.Cells.Item(44, 8).Value = 645
$ws.Cells.Item(71, 2).Value = 9573
$ws.Cells.Item(71, 3).Value = 316
$ws.Cells.Item(71, 4).Value = 4606
$ws.Cells.Item(71, 5).Value = 4365
$ws.Cells.Item(71, 7).Value = 30
$ws.Cells.Item(71, 8).Value = 602
$ws.Cells.Item(87, 1).Value = "Gabon"
$ws.Cells.Item(87, 2).Value = 5513
$ws.Cells.Item(87, 3).Value = 119
$ws.Cells.Item(87, 4).Value = 2508
$ws.Cells.Item(87, 5).Value = 2963
$ws.Cells.Item(87, 8).Value = 42
$ws.Cells.Item(88, 1).Value = "Guinea"
$ws.Cells.Item(88, 2).Value = 5404
$ws.Cells.Item(88, 3).Value = 13
$ws.Cells.Item(88, 4).Value = 4346
$ws.Cells.Item(88, 5).Value = 1025
$ws.Cells.Item(88, 8).Value = 33
$ws.Cells.Item(124, 1).Value = "Congo"
$ws.Cells.Item(124, 2).Value = 1382
$ws.Cells.Item(124, 3).Value = 295
$ws.Cells.Item(124, 4).Value = 486
$ws.Cells.Item(124, 5).Value = 855
$ws.Cells.Item(124, 7).Value = 4
$ws.Cells.Item(124, 8).Value = 41
$ws.Cells.Item(125, 1).Value = "Cabo Verde"
$ws.Cells.Item(125, 2).Value = 1267
$ws.Cells.Item(125, 3).Value = 40
$ws.Cells.Item(125, 4).Value = 629
$ws.Cells.Item(125, 5).Value = 623
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 15
$ws.Cells.Item(126, 1).Value = "Malaui"
$ws.Cells.Item(126, 2).Value = 1265
$ws.Cells.Item(126, 3).Value = 41
$ws.Cells.Item(126, 4).Value = 260
$ws.Cells.Item(126, 5).Value = 989
$ws.Cells.Item(126, 7).Value = 2
$ws.Cells.Item(126, 8).Value = 16
$ws.Cells.Item(127, 1).Value = "Hong Kong"
$ws.Cells.Item(127, 2).Value = 1234
$ws.Cells.Item(127, 3).Value = 28
$ws.Cells.Item(127, 4).Value = 1117
$ws.Cells.Item(127, 5).Value = 110
$ws.Cells.Item(127, 8).Value = 7
$ws.Cells.Item(128, 1).Value = "Benin"
$ws.Cells.Item(128, 2).Value = 1199
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(128, 4).Value = 333
$ws.Cells.Item(128, 5).Value = 845
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 21
$ws.Cells.Item(129, 1).Value = "Yemen"
$ws.Cells.Item(129, 2).Value = 1190
$ws.Cells.Item(129, 3).Value = 32
$ws.Cells.Item(129, 4).Value = 504
$ws.Cells.Item(129, 5).Value = 368
$ws.Cells.Item(129, 7).Value = 6
$ws.Cells.Item(129, 8).Value = 318
$ws.Cells.Item(130, 1).Value = "Tunez"
$ws.Cells.Item(130, 2).Value = 1175
$ws.Cells.Item(130, 4).Value = 1038
$ws.Cells.Item(130, 5).Value = 87
$ws.Cells.Item(130, 8).Value = 50
$ws.Cells.Item(131, 1).Value = "Jordania"
$ws.Cells.Item(131, 2).Value = 1133
$ws.Cells.Item(131, 3).Value = 1
$ws.Cells.Item(131, 4).Value = 886
$ws.Cells.Item(131, 5).Value = 238
$ws.Cells.Item(131, 8).Value = 9
$ws.Cells.Item(132, 1).Value = "Letonia"
$ws.Cells.Item(132, 2).Value = 1121
$ws.Cells.Item(132, 3).Value = 3
$ws.Cells.Item(132, 4).Value = 974
$ws.Cells.Item(132, 5).Value = 117
$ws.Cells.Item(132, 8).Value = 30
$ws.Cells.Item(141, 1).Value = "Libia"
$ws.Cells.Item(141, 2).Value = 874
$ws.Cells.Item(141, 3).Value = 50
$ws.Cells.Item(141, 4).Value = 223
$ws.Cells.Item(141, 5).Value = 626
$ws.Cells.Item(141, 7).Value = 1
$ws.Cells.Item(141, 8).Value = 25
$ws.Cells.Item(142, 1).Value = "Republica del Chad"
$ws.Cells.Item(142, 2).Value = 866
$ws.Cells.Item(142, 4).Value = 785
$ws.Cells.Item(142, 5).Value = 7
$ws.Cells.Item(142, 8).Value = 74
$ws.Cells.Item(143, 1).Value = "Principado de Andorra"
$ws.Cells.Item(143, 2).Value = 855
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 799
$ws.Cells.Item(143, 5).Value = 4
$ws.Cells.Item(143, 8).Value = 52
$ws.Cells.Item(144, 1).Value = "Suazilandia"
$ws.Cells.Item(144, 2).Value = 840
$ws.Cells.Item(144, 3).Value = 28
$ws.Cells.Item(144, 4).Value = 418
$ws.Cells.Item(144, 5).Value = 411
$ws.Cells.Item(144, 8).Value = 11
$ws.Cells.Item(151, 2).Value = 661
$ws.Cells.Item(151, 3).Value = 11
$ws.Cells.Item(151, 4).Value = 414
$ws.Cells.Item(151, 5).Value = 233
$ws.Cells.Item(152, 2).Value = 605
$ws.Cells.Item(152, 3).Value = 14
$ws.Cells.Item(152, 4).Value = 166
$ws.Cells.Item(152, 5).Value = 432
$ws.Cells.Item(166, 5).Value = 119
$ws.Cells.Item(166, 7).Value = 1
$ws.Cells.Item(166, 8).Value = 13
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"
